$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newRows = @(
    @(44304, 4, 43, 178.2605090788492),
    @(44305, 5, 43, 178.2605090788492),
    @(44306, 2, 39, 161.6781361412818),
    @(44307, 0, 38, 157.53254290689)
)

$lastRow = 229
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $lastRow + $i + 1
    $srcRow = $lastRow

    # Copy formatting/style from the previous row so the new row matches existing styling
    $ws.Range("A$srcRow`:D$srcRow").Copy($ws.Range("A$targetRow`:D$targetRow"))

    $vals = $newRows[$i]
    $ws.Cells.Item($targetRow, 1).Value = $vals[0]
    $ws.Cells.Item($targetRow, 2).Value = $vals[1]
    $ws.Cells.Item($targetRow, 3).Value = $vals[2]
    $ws.Cells.Item($targetRow, 4).Value = $vals[3]
}
